$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row correct-answer total (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (B12): 51 -> 85, corresponding marks fraction (E12): 47/84 -> 85/140
$ws.Range("B12").Value = 85
$ws.Range("E12").Value = "85/140"
